$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 331, shifting all existing data (rows 331-383)
# down by one (new last row becomes 384).
$ws.Rows.Item(331).Insert()

# Populate the newly inserted row 331 with the new record.
$ws.Range("A331").Value = 5
$ws.Range("B331").Value = "Macroferia Regional de Talca"
$ws.Range("C331").Value = "Maule"
$ws.Range("D331").Value = 44474
$ws.Range("E331").Value = 7
$ws.Range("F331").Value = 100112004
$ws.Range("G331").Value = "Cebolla"
$ws.Range("H331").Value = "Sin especificar"
$ws.Range("I331").Value = "1a nueva(o)"
$ws.Range("J331").Value = 50000
$ws.Range("K331").Value = 1200
$ws.Range("L331").Value = 1200
$ws.Range("M331").Value = 1200
$ws.Range("N331").Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O331").Value = "Región de O'Higgins"
$ws.Range("P331").Value = 120
$ws.Range("Q331").Value = 10
$ws.Range("R331").Value = "Hortaliza"
